$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value renders as a plain number and must be forced to
# Text format first so Excel does not coerce it to a numeric value (which
# would destroy formatting such as trailing zeros, e.g. "39.60" -> 39.6).
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '588.19'
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '150.34'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '5.71'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '27.57'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '12.21'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '344.41'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '67.26'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '9.23'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.67'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '558.44'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.23'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '166.64'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '166.14'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '39.60'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '22.81'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.628'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0960'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '19.20'

# Remaining cells (percentages and dotted/odd-format prices) are already
# non-numeric text, so they can be assigned directly.
$ws.Range('D2').Value = '63.311.21'
$ws.Range('E2').Value = '  -0.87%  '
$ws.Range('D3').Value = '2.582.75'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('E5').Value = '  -3.23%  '
$ws.Range('E6').Value = '  +2.17%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E9').Value = '  +0.68%  '
$ws.Range('E10').Value = '  +2.01%  '
$ws.Range('E11').Value = '  -0.10%  '
$ws.Range('E12').Value = '  -0.38%  '
$ws.Range('E13').Value = '  +0.38%  '
$ws.Range('D14').Value = '3.045.09'
$ws.Range('E14').Value = '  -2.17%  '
$ws.Range('D15').Value = '63.118.81'
$ws.Range('E15').Value = '  -0.86%  '
$ws.Range('E16').Value = '  +5.17%  '
$ws.Range('D17').Value = '2.588.31'
$ws.Range('E17').Value = '  -2.52%  '
$ws.Range('E18').Value = '  +3.78%  '
$ws.Range('E19').Value = '  +3.56%  '
$ws.Range('E20').Value = '  -0.73%  '
$ws.Range('E21').Value = '  -1.07%  '
$ws.Range('E22').Value = '  -0.04%  '
$ws.Range('E23').Value = '  +1.42%  '
$ws.Range('E24').Value = '  +0.91%  '
$ws.Range('E25').Value = '  -0.32%  '
$ws.Range('E26').Value = '  -1.36%  '
$ws.Range('E27').Value = '  -1.06%  '
$ws.Range('E28').Value = '  -0.93%  '
$ws.Range('E29').Value = '  +1.04%  '
$ws.Range('E30').Value = '  -0.01%  '
$ws.Range('E31').Value = '  -1.04%  '
$ws.Range('D32').Value = '0.0₃0856'
$ws.Range('E32').Value = '  +0.44%  '
$ws.Range('E33').Value = '  -0.97%  '
$ws.Range('E34').Value = '  -1.08%  '
$ws.Range('E35').Value = '  -1.31%  '
$ws.Range('E36').Value = '  +1.77%  '
$ws.Range('E37').Value = '  +0.05%  '
$ws.Range('E38').Value = '  +1.74%  '
$ws.Range('E39').Value = '  -1.74%  '
$ws.Range('E41').Value = '  +0.58%  '
$ws.Range('E43').Value = '  +5.33%  '
$ws.Range('E44').Value = '  +3.87%  '
$ws.Range('E45').Value = '  +2.55%  '
$ws.Range('E46').Value = '  +5.95%  '
$ws.Range('E47').Value = '  +0.03%  '
$ws.Range('E48').Value = '  +1.94%  '
$ws.Range('E49').Value = '  +0.21%  '
$ws.Range('E50').Value = '  +1.95%  '
$ws.Range('E51').Value = '  +18.33%  '
